$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, pushing existing rows 46:51 down to 47:52
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the new weekly data point
$ws.Range("A46").Value = 8
$ws.Range("B46").Value = "Terminal La Palmera de La Serena"
$ws.Range("C46").Value = "Coquimbo"
$ws.Range("D46").Value = 44505
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 100112052
$ws.Range("G46").Value = "Albahaca"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = 3750
$ws.Range("N46").Value = "$/paquete"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 3750
$ws.Range("Q46").Value = 1
$ws.Range("R46").Value = "Hortaliza"
